# Tag-3_3-Abschluss.pptx edit:
#  - the footer "author" text box (inherited from the slide layout) changes
#    from "Daniel Krämer" to "Patrick Möbius"
#  - the footer date field (inherited from the slide master) changes its
#    cached text from "15.03.2019" to "14.05.2021"

$p = $ppt.ActivePresentation

# --- 1) Author name, lives on the slide layout shared by the slides ---
$layout = $p.Slides.Item(1).CustomLayout
$authorShape = $layout.Shapes.Item("Text Box 24")
$tf = $authorShape.TextFrame
$tr = $tf.TextRange

# This textbox has two paragraphs: the address line, then the author line.
# Only touch the author paragraph (the last one) so the address text and
# its runs are left completely untouched.
$paraCount = $tr.Paragraphs().Count
$authorPara = $tr.Paragraphs($paraCount, 1)
$authorPara.Text = "Patrick Möbius"

# --- 2) Footer date field, lives on the slide master ---
$master = $p.SlideMaster
$dateShape = $master.Shapes.Item("Rectangle 6")
$dateShape.TextFrame.TextRange.Text = "14.05.2021"
